$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, shifting existing rows 29:122 down to 30:123.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record.
$ws.Range("A29").Value = 7
$ws.Range("B29").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C29").Value = "Ñuble"
$ws.Range("D29").Value = 44980
$ws.Range("E29").Value = 16
$ws.Range("F29").Value = 100112031
$ws.Range("G29").Value = "Poroto verde"
$ws.Range("H29").Value = "Magnum"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 50
$ws.Range("K29").Value = 25000
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 25000
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Provincia de Diguillín"
$ws.Range("P29").Value = 1000
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
